$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("G1:H1").EntireColumn.Insert()
$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"
$ws.Range("F2").Value = 0.0700106
$ws.Range("G2").Value = 0.0365131
$ws.Range("H2").Value = 0.106524
$ws.Range("I2").Value = 10.8
$ws.Range("F3").Value = 0.0698885
$ws.Range("G3").Value = 0.0371189
$ws.Range("H3").Value = 0.107007
$ws.Range("I3").Value = 10.8
$ws.Range("F4").Value = 0.0710265
$ws.Range("G4").Value = 0.0383534
$ws.Range("H4").Value = 0.10938
$ws.Range("I4").Value = 10.8
$ws.Range("F5").Value = 0.0694585
$ws.Range("G5").Value = 0.036635
$ws.Range("H5").Value = 0.106094
$ws.Range("I5").Value = 10.7
$ws.Range("F6").Value = 0.0702814
$ws.Range("G6").Value = 0.037701
$ws.Range("H6").Value = 0.107982
$ws.Range("I6").Value = 10.8
$ws.Range("F7").Value = 0.0700087
$ws.Range("G7").Value = 0.0366247
$ws.Range("H7").Value = 0.106633
$ws.Range("I7").Value = 10.8
$ws.Range("F8").Value = 0.0695495
$ws.Range("G8").Value = 0.0365601
$ws.Range("H8").Value = 0.10611
$ws.Range("I8").Value = 10.7
$ws.Range("F9").Value = 0.0698445
$ws.Range("G9").Value = 0.037205
$ws.Range("H9").Value = 0.10705
$ws.Range("I9").Value = 10.8
$ws.Range("F10").Value = 0.0710091
$ws.Range("G10").Value = 0.0374241
$ws.Range("H10").Value = 0.108433
$ws.Range("I10").Value = 10.7
$ws.Range("F11").Value = 0.0696036
$ws.Range("G11").Value = 0.036656
$ws.Range("H11").Value = 0.10626
$ws.Range("I11").Value = 10.7
$ws.Range("F12").Value = 0.07006809
$ws.Range("G12").Value = 0.03707913
$ws.Range("H12").Value = 0.1071473
$ws.Range("I12").Value = 10.76
$ws.Range("J2:K12").ClearContents()

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("G1:H1").EntireColumn.Insert()
$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"
$ws.Range("F2").Value = 0.0987603
$ws.Range("G2").Value = 0.0452667
$ws.Range("H2").Value = 0.144027
$ws.Range("I2").Value = 16.8
$ws.Range("F3").Value = 0.100749
$ws.Range("G3").Value = 0.0449665
$ws.Range("H3").Value = 0.145715
$ws.Range("I3").Value = 16.7
$ws.Range("F4").Value = 0.0971588
$ws.Range("G4").Value = 0.0440111
$ws.Range("H4").Value = 0.14117
$ws.Range("I4").Value = 16.7
$ws.Range("F5").Value = 0.0980915
$ws.Range("G5").Value = 0.0440222
$ws.Range("H5").Value = 0.142114
$ws.Range("I5").Value = 16.8
$ws.Range("F6").Value = 0.0989569
$ws.Range("G6").Value = 0.0448444
$ws.Range("H6").Value = 0.143801
$ws.Range("I6").Value = 16.7
$ws.Range("F7").Value = 0.0994543
$ws.Range("G7").Value = 0.0451059
$ws.Range("H7").Value = 0.14456
$ws.Range("I7").Value = 16.7
$ws.Range("F8").Value = 0.0992628
$ws.Range("G8").Value = 0.0448333
$ws.Range("H8").Value = 0.144096
$ws.Range("I8").Value = 16.7
$ws.Range("F9").Value = 0.0979723
$ws.Range("G9").Value = 0.04573
$ws.Range("H9").Value = 0.143702
$ws.Range("I9").Value = 16.7
$ws.Range("F10").Value = 0.0974166
$ws.Range("G10").Value = 0.043687
$ws.Range("H10").Value = 0.141104
$ws.Range("I10").Value = 16.7
$ws.Range("F11").Value = 0.0973398
$ws.Range("G11").Value = 0.0444539
$ws.Range("H11").Value = 0.141794
$ws.Range("I11").Value = 16.7
$ws.Range("F12").Value = 0.09851623
$ws.Range("G12").Value = 0.0446921
$ws.Range("H12").Value = 0.1432083
$ws.Range("I12").Value = 16.72
$ws.Range("J2:K12").ClearContents()

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("G1:H1").EntireColumn.Insert()
$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"
$ws.Range("F2").Value = 0.126339
$ws.Range("G2").Value = 0.0481958
$ws.Range("H2").Value = 0.174535
$ws.Range("I2").Value = 24.6
$ws.Range("F3").Value = 0.12769
$ws.Range("G3").Value = 0.0486617
$ws.Range("H3").Value = 0.176351
$ws.Range("I3").Value = 24.2
$ws.Range("F4").Value = 0.129951
$ws.Range("G4").Value = 0.0505224
$ws.Range("H4").Value = 0.180473
$ws.Range("I4").Value = 24.6
$ws.Range("F5").Value = 0.133103
$ws.Range("G5").Value = 0.0477876
$ws.Range("H5").Value = 0.18089
$ws.Range("I5").Value = 24.2
$ws.Range("F6").Value = 0.131039
$ws.Range("G6").Value = 0.0496911
$ws.Range("H6").Value = 0.18073
$ws.Range("I6").Value = 24.8
$ws.Range("F7").Value = 0.135257
$ws.Range("G7").Value = 0.0510407
$ws.Range("H7").Value = 0.186298
$ws.Range("I7").Value = 24.2
$ws.Range("F8").Value = 0.128737
$ws.Range("G8").Value = 0.049667
$ws.Range("H8").Value = 0.178404
$ws.Range("I8").Value = 24.6
$ws.Range("F9").Value = 0.12783
$ws.Range("G9").Value = 0.0498899
$ws.Range("H9").Value = 0.17772
$ws.Range("I9").Value = 24.6
$ws.Range("F10").Value = 0.130924
$ws.Range("G10").Value = 0.0501136
$ws.Range("H10").Value = 0.181038
$ws.Range("I10").Value = 24.5
$ws.Range("F11").Value = 0.133274
$ws.Range("G11").Value = 0.0510449
$ws.Range("H11").Value = 0.184319
$ws.Range("I11").Value = 24.2
$ws.Range("F12").Value = 0.1304144
$ws.Range("G12").Value = 0.04966147
$ws.Range("H12").Value = 0.1800758
$ws.Range("I12").Value = 24.45
$ws.Range("J2:K12").ClearContents()

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("G1:H1").EntireColumn.Insert()
$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"
$ws.Range("F2").Value = 0.159301
$ws.Range("G2").Value = 0.052887
$ws.Range("H2").Value = 0.212188
$ws.Range("I2").Value = 29.3
$ws.Range("F3").Value = 0.155994
$ws.Range("G3").Value = 0.0546329
$ws.Range("H3").Value = 0.210627
$ws.Range("I3").Value = 29.3
$ws.Range("F4").Value = 0.159755
$ws.Range("G4").Value = 0.0522819
$ws.Range("H4").Value = 0.212037
$ws.Range("I4").Value = 29.2
$ws.Range("F5").Value = 0.161734
$ws.Range("G5").Value = 0.0543456
$ws.Range("H5").Value = 0.216079
$ws.Range("I5").Value = 29.5
$ws.Range("F6").Value = 0.154287
$ws.Range("G6").Value = 0.0511792
$ws.Range("H6").Value = 0.205467
$ws.Range("I6").Value = 29
$ws.Range("F7").Value = 0.155964
$ws.Range("G7").Value = 0.0519663
$ws.Range("H7").Value = 0.20793
$ws.Range("I7").Value = 29.7
$ws.Range("F8").Value = 0.163388
$ws.Range("G8").Value = 0.0524677
$ws.Range("H8").Value = 0.215855
$ws.Range("I8").Value = 29.8
$ws.Range("F9").Value = 0.163394
$ws.Range("G9").Value = 0.0556584
$ws.Range("H9").Value = 0.219052
$ws.Range("I9").Value = 29.7
$ws.Range("F10").Value = 0.156548
$ws.Range("G10").Value = 0.0537551
$ws.Range("H10").Value = 0.210303
$ws.Range("I10").Value = 29.3
$ws.Range("F11").Value = 0.158463
$ws.Range("G11").Value = 0.0547844
$ws.Range("H11").Value = 0.213247
$ws.Range("I11").Value = 29.8
$ws.Range("F12").Value = 0.1588828
$ws.Range("G12").Value = 0.05339585
$ws.Range("H12").Value = 0.2122785
$ws.Range("I12").Value = 29.46
$ws.Range("J2:K12").ClearContents()
